$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps being stored as text (it already holds
# values like "42.654.07" that are not valid numbers), so that values such
# as "24.60" or "18.10" are not silently reinterpreted as numbers and lose
# their trailing zero / multi-dot formatting.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "42.680.79"
$ws.Range("E2").Value = "  -0.52%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.294.27"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "301.87"

# Row 6 - Solana
$ws.Range("D6").Value = "95.97"
$ws.Range("E6").Value = "  -1.58%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.96%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -2.04%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "34.75"
$ws.Range("E10").Value = "  -2.88%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").Value = "  -0.80%  "

# Row 12 - Chainlink
$ws.Range("D12").Value = "18.56"
$ws.Range("E12").Value = "  +4.57%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +2.28%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +1.00%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.651.81"
$ws.Range("E15").Value = "  -0.32%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.299.59"
$ws.Range("E16").Value = "  -0.05%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "0.773"
$ws.Range("E17").Value = "  -0.76%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "42.595.12"
$ws.Range("E18").Value = "  -0.67%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").Value = "12.77"
$ws.Range("E19").Value = "  +1.26%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -1.81%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -1.71%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "67.02"
$ws.Range("E22").Value = "  -1.35%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "235.82"
$ws.Range("E23").Value = "  -2.55%  "

# Row 24 - ImmutableX
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  -0.38%  "

# Row 26 - PancakeSwap
$ws.Range("D26").Value = "2.38"
$ws.Range("E26").Value = "  -1.71%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "24.60"
$ws.Range("E27").Value = "  -1.64%  "

# Row 28 - Monero
$ws.Range("D28").Value = "167.27"
$ws.Range("E28").Value = "  +0.63%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  +0.60%  "

# Row 30 - Cosmos
$ws.Range("D30").Value = "8.98"
$ws.Range("E30").Value = "  -0.47%  "

# Row 31 - InjectiveProtocol
$ws.Range("E31").Value = "  +0.41%  "

# Row 32 - FirstDigitalUSD
$ws.Range("E32").Value = "  +0.00%  "

# Row 33 - Celestia
$ws.Range("D33").Value = "17.73"
$ws.Range("E33").Value = "  +1.76%  "

# Row 34 - Filecoin
$ws.Range("D34").Value = "4.95"
$ws.Range("E34").Value = "  -0.97%  "

# Row 35 - RenderToken
$ws.Range("E35").Value = "  -6.56%  "

# Row 36 - WEMIXToken
$ws.Range("D36").Value = "2.35"
$ws.Range("E36").Value = "  -1.80%  "

# Row 37 - Hedera
$ws.Range("E37").Value = "  -0.10%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  -0.49%  "

# Row 39 - ARBITRUM
$ws.Range("D39").Value = "1.74"
$ws.Range("E39").Value = "  -1.55%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  -1.10%  "

# Row 41 - LidoDAOToken
$ws.Range("E41").Value = "  -2.67%  "

# Row 42 - Maker
$ws.Range("D42").Value = "1.992.09"
$ws.Range("E42").Value = "  -0.53%  "

# Row 43 - VeChain
$ws.Range("E43").Value = "  -2.00%  "

# Row 44 - FraxShare
$ws.Range("D44").Value = "10.18"
$ws.Range("E44").Value = "  +0.46%  "

# Row 45 & 46 - EnergySwap and ApeXProtocol swap places with updated data
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "2.14"
$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "18.10"
$ws.Range("E46").Value = "  +5.54%  "

# Row 47 - NEARProtocol
$ws.Range("E47").Value = "  -0.60%  "

# Row 48 - MultiversX
$ws.Range("D48").Value = "53.35"
$ws.Range("E48").Value = "  +0.23%  "

# Row 49 - HuobiToken
$ws.Range("E49").Value = "  +4.59%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "2.518.70"
$ws.Range("E50").Value = "  -0.28%  "

# Row 51 - BitcoinSV
$ws.Range("D51").Value = "70.70"
$ws.Range("E51").Value = "  -1.99%  "
